$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 6783
$ws.Range("J16").Value = 6783
$ws.Range("L16").Value = 6783
$ws.Range("N16").Value = -7243
$ws.Range("H87").Value = 40000
$ws.Range("J87").Value = 40000
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42496
$ws.Range("H90").Value = 40000
$ws.Range("J90").Value = 40000
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -132480
$ws.Range("H138").Value = 4690.636
$ws.Range("J138").Value = 5066.6665
$ws.Range("L138").Value = 15199.9995
$ws.Range("N138").Value = -25479.9995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 21997.5
$ws.Range("J37").Value = 23330.555
$ws.Range("L37").Value = 23330.555
$ws.Range("N37").Value = -23876.555
$ws.Range("H45").Value = 2002
$ws.Range("I45").Value = 2002
$ws.Range("K45").Value = 2002
$ws.Range("M45").Value = -1625
$ws.Range("H63").Value = 3473.5454
$ws.Range("I63").Value = 2802.25
$ws.Range("J63").Value = 3857.1428
$ws.Range("K63").Value = 2802.25
$ws.Range("L63").Value = 3857.1428
$ws.Range("M63").Value = -2116.25
$ws.Range("N63").Value = -5229.1428
$ws.Range("H66").Value = 3473.5454
$ws.Range("I66").Value = 2802.25
$ws.Range("J66").Value = 3857.1428
$ws.Range("K66").Value = 14011.25
$ws.Range("L66").Value = 19285.714
$ws.Range("M66").Value = -10579.25
$ws.Range("N66").Value = -26149.714
$ws.Range("H74").Value = 2273.9473
$ws.Range("I74").Value = 2222.5
$ws.Range("J74").Value = 3200
$ws.Range("K74").Value = 2222.5
$ws.Range("L74").Value = 3200
$ws.Range("M74").Value = -1348.5
$ws.Range("N74").Value = -4948
$ws.Range("H77").Value = 2273.9473
$ws.Range("I77").Value = 2222.5
$ws.Range("J77").Value = 3200
$ws.Range("K77").Value = 11112.5
$ws.Range("L77").Value = 16000
$ws.Range("M77").Value = -6744.5
$ws.Range("N77").Value = -24736
$ws.Range("H122").Value = 500
$ws.Range("J122").Value = 500
$ws.Range("L122").Value = 1500
$ws.Range("N122").Value = -6400

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 14996.875
$ws.Range("J35").Value = 14996.875
$ws.Range("L35").Value = 14996.875
$ws.Range("N35").Value = -15616.875
$ws.Range("H82").Value = 27399.625
$ws.Range("J82").Value = 39997.5
$ws.Range("L82").Value = 39997.5
$ws.Range("N82").Value = -40763.5
$ws.Range("H85").Value = 27399.625
$ws.Range("J85").Value = 39997.5
$ws.Range("L85").Value = 39997.5
$ws.Range("N85").Value = -42649.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9949.799999999999
$ws.Range("I16").Value = 8249.666999999999
$ws.Range("J16").Value = 12500
$ws.Range("K16").Value = 8249.666999999999
$ws.Range("L16").Value = 12500
$ws.Range("M16").Value = -7962.666999999999
$ws.Range("N16").Value = -13074
$ws.Range("H22").Value = 766.6667
$ws.Range("I22").Value = 740
$ws.Range("K22").Value = 740
$ws.Range("M22").Value = -390
$ws.Range("H33").Value = 2015.25
$ws.Range("I33").Value = 2015.25
$ws.Range("K33").Value = 2015.25
$ws.Range("M33").Value = -1636.25
$ws.Range("H41").Value = 19085
$ws.Range("J41").Value = 20859.445
$ws.Range("L41").Value = 20859.445
$ws.Range("N41").Value = -21715.445
$ws.Range("H86").Value = 2786.25
$ws.Range("I86").Value = 2614.2856
$ws.Range("J86").Value = 3990
$ws.Range("K86").Value = 2614.2856
$ws.Range("L86").Value = 3990
$ws.Range("M86").Value = -1491.2856
$ws.Range("N86").Value = -6236
$ws.Range("H89").Value = 2786.25
$ws.Range("I89").Value = 2614.2856
$ws.Range("J89").Value = 3990
$ws.Range("K89").Value = 13071.428
$ws.Range("L89").Value = 19950
$ws.Range("M89").Value = -7455.428
$ws.Range("N89").Value = -31182
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0
$ws.Range("H113").Value = 9949.799999999999
$ws.Range("I113").Value = 8249.666999999999
$ws.Range("J113").Value = 12500
$ws.Range("K113").Value = 8249.666999999999
$ws.Range("L113").Value = 12500
$ws.Range("M113").Value = -6079.666999999999
$ws.Range("N113").Value = -16840
$ws.Range("H122").Value = 3941.6667
$ws.Range("I122").Value = 4029.4707
$ws.Range("K122").Value = 12088.4121
$ws.Range("M122").Value = -9638.4121
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 0

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 501.73685
$ws.Range("I5").Value = 512.1667
$ws.Range("J5").Value = 314
$ws.Range("K5").Value = 1536.5001
$ws.Range("L5").Value = 942
$ws.Range("M5").Value = -1424.5001
$ws.Range("N5").Value = -1166
$ws.Range("H80").Value = 751.5
$ws.Range("I80").Value = 500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -564
$ws.Range("H83").Value = 751.5
$ws.Range("I83").Value = 500
$ws.Range("K83").Value = 4500
$ws.Range("M83").Value = 180
$ws.Range("H109").Value = 1942.3636
$ws.Range("J109").Value = 2796.2
$ws.Range("L109").Value = 8388.599999999999
$ws.Range("N109").Value = -10468.6
$ws.Range("H122").Value = 1853.5714
$ws.Range("I122").Value = 1480
$ws.Range("K122").Value = 13320
$ws.Range("M122").Value = -10870
$ws.Range("H131").Value = 20315.166
$ws.Range("J131").Value = 4365.625
$ws.Range("L131").Value = 13096.875
$ws.Range("N131").Value = -23176.875
$ws.Range("H132").Value = 801.75
$ws.Range("I132").Value = 801.75
$ws.Range("K132").Value = 7215.75
$ws.Range("M132").Value = -4685.75
$ws.Range("H135").Value = 501.73685
$ws.Range("I135").Value = 512.1667
$ws.Range("J135").Value = 314
$ws.Range("K135").Value = 4609.5003
$ws.Range("L135").Value = 2826
$ws.Range("M135").Value = -2074.5003
$ws.Range("N135").Value = -7896

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 15428.571
$ws.Range("I46").Value = 8000
$ws.Range("J46").Value = 18400
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 18400
$ws.Range("M46").Value = -7844
$ws.Range("N46").Value = -18712

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H46").Value = 3999.3333
$ws.Range("J46").Value = 3998
$ws.Range("L46").Value = 3998
$ws.Range("N46").Value = -4374
$ws.Range("H76").Value = 2800.8
$ws.Range("J76").Value = 2800.8
$ws.Range("L76").Value = 2800.8
$ws.Range("N76").Value = -3476.8
$ws.Range("H79").Value = 2800.8
$ws.Range("J79").Value = 2800.8
$ws.Range("L79").Value = 2800.8
$ws.Range("N79").Value = -5140.8
$ws.Range("H93").Value = 1672.5
$ws.Range("I93").Value = 1672.5
$ws.Range("K93").Value = 1672.5
$ws.Range("M93").Value = -424.5
$ws.Range("H122").Value = 4337
$ws.Range("I122").Value = 4275
$ws.Range("K122").Value = 12825
$ws.Range("M122").Value = -10375

Write-Host "Applied 199 cell changes across 7 sheets"